# Update "Ciudades" worksheet with refreshed COVID case counts and timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1.
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 20:05"

# Row 4 - Madrid
$ws.Range("C4").Value = 40736
$ws.Range("D4").Value = 16703

# Row 5 - Cataluña
$ws.Range("C5").Value = 26203
$ws.Range("D5").Value = 23762

# Row 6 - Castilla y Leon
$ws.Range("D6").Value = 7843

# Row 7 - Castilla-La Mancha
$ws.Range("C7").Value = 6392
$ws.Range("D7").Value = 7358

# Row 9 - Andalucia
$ws.Range("C9").Value = 10671
$ws.Range("D9").Value = 429

# Row 14 - Aragon
$ws.Range("C14").Value = 3772
$ws.Range("D14").Value = 893

# Row 16 - Navarra
$ws.Range("C16").Value = 3905
$ws.Range("D16").Value = 770

# Row 20 - Salamanca
$ws.Range("C20").Value = 3107
$ws.Range("D20").Value = 566

# Row 32 - Granada
$ws.Range("C32").Value = 1063
$ws.Range("D32").Value = 989

# Row 33 - Asturias
$ws.Range("C33").Value = 1537
$ws.Range("D33").Value = 604
